$d = $word.ActiveDocument

$pairs = @(
    @("491÷7=70, 1", "476÷4=119, 0"),
    @("190÷9=21, 1", "231÷8=28, 7"),
    @("799÷4=199, 3", "464÷3=154, 2"),
    @("113÷5=22, 3", "750÷5=150, 0"),
    @("165÷4=41, 1", "314÷4=78, 2"),
    @("930÷3=310, 0", "513÷8=64, 1"),
    @("582÷2=291, 0", "481÷6=80, 1"),
    @("530÷2=265, 0", "428÷8=53, 4"),
    @("633÷6=105, 3", "930÷7=132, 6"),
    @("306÷5=61, 1", "781÷8=97, 5"),
    @("826÷9=91, 7", "653÷5=130, 3"),
    @("480÷2=240, 0", "476÷8=59, 4"),
    @("347÷4=86, 3", "938÷2=469, 0"),
    @("271÷3=90, 1", "360÷2=180, 0"),
    @("204÷6=34, 0", "892÷9=99, 1"),
    @("741÷7=105, 6", "528÷2=264, 0"),
    @("654÷6=109, 0", "551÷2=275, 1"),
    @("490÷6=81, 4", "952÷4=238, 0"),
    @("660÷8=82, 4", "152÷2=76, 0"),
    @("604÷3=201, 1", "382÷4=95, 2"),
    @("620÷4=155, 0", "505÷9=56, 1"),
    @("646÷8=80, 6", "785÷3=261, 2"),
    @("123÷4=30, 3", "765÷3=255, 0"),
    @("574÷6=95, 4", "199÷9=22, 1"),
    @("331÷3=110, 1", "684÷5=136, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
